$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Add a new row of data (row 4) to Sheet1, reusing existing strings str7/str8/str9
$ws1.Range("F4").Value = "str7"
$ws1.Range("G4").Value = "str8"
$ws1.Range("H4").Value = "str9"

# Update selection on Sheet1 (no longer the tab-selected sheet)
$ws1.Range("H14").Select()

# Update selection on Sheet2 and make it the active/selected sheet
$ws2.Activate()
$ws2.Range("C1").Select()
